$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 47.86240033333333
$ws.Cells.Item(2, 8).Value2 = 143.587201
$ws.Cells.Item(2, 9).Value2 = 0.1228118231805696
$ws.Cells.Item(2, 10).Value2 = 0.1228118231805696
$ws.Cells.Item(2, 13).Value2 = 42.99144133333333
$ws.Cells.Item(2, 14).Value2 = 128.974324
$ws.Cells.Item(2, 15).Value2 = 0.2509605789056467
$ws.Cells.Item(2, 16).Value2 = 0.2509605789056467
$ws.Cells.Item(2, 17).Value2 = 2057.673576003014
$ws.Cells.Item(2, 18).Value2 = 18519.06218402712
$ws.Cells.Item(2, 19).Value2 = 0.03082092624185367
$ws.Cells.Item(2, 20).Value2 = 0.03082092624185367
$ws.Cells.Item(3, 7).Value2 = 47.86240033333333
$ws.Cells.Item(3, 8).Value2 = 143.587201
$ws.Cells.Item(3, 9).Value2 = 0.1228118231805696
$ws.Cells.Item(3, 10).Value2 = 0.1228118231805696
$ws.Cells.Item(3, 15).Value2 = 0.3246035785195009
$ws.Cells.Item(3, 16).Value2 = 0.324603578519501
$ws.Cells.Item(3, 17).Value2 = 2661.486553418881
$ws.Cells.Item(3, 18).Value2 = 23953.37898076992
$ws.Cells.Item(3, 19).Value2 = 0.0398651572889171
$ws.Cells.Item(3, 20).Value2 = 0.0398651572889171
$ws.Cells.Item(4, 7).Value2 = 47.86240033333333
$ws.Cells.Item(4, 8).Value2 = 143.587201
$ws.Cells.Item(4, 9).Value2 = 0.1228118231805696
$ws.Cells.Item(4, 10).Value2 = 0.1228118231805696
$ws.Cells.Item(4, 13).Value2 = 36.72715
$ws.Cells.Item(4, 14).Value2 = 110.18145
$ws.Cells.Item(4, 15).Value2 = 0.2143930638214748
$ws.Cells.Item(4, 16).Value2 = 0.2143930638214748
$ws.Cells.Item(4, 17).Value2 = 1757.849556402383
$ws.Cells.Item(4, 18).Value2 = 15820.64600762145
$ws.Cells.Item(4, 19).Value2 = 0.02633000304518354
$ws.Cells.Item(4, 20).Value2 = 0.02633000304518354
$ws.Cells.Item(5, 7).Value2 = 47.86240033333333
$ws.Cells.Item(5, 8).Value2 = 143.587201
$ws.Cells.Item(5, 9).Value2 = 0.1228118231805696
$ws.Cells.Item(5, 10).Value2 = 0.1228118231805696
$ws.Cells.Item(5, 13).Value2 = 7.381512666666667
$ws.Cells.Item(5, 14).Value2 = 22.144538
$ws.Cells.Item(5, 15).Value2 = 0.04308924368603855
$ws.Cells.Item(5, 16).Value2 = 0.04308924368603857
$ws.Cells.Item(5, 17).Value2 = 353.2969143175709
$ws.Cells.Item(5, 18).Value2 = 3179.672228858138
$ws.Cells.Item(5, 19).Value2 = 0.005291868576554243
$ws.Cells.Item(5, 20).Value2 = 0.005291868576554244
$ws.Cells.Item(6, 7).Value2 = 47.86240033333333
$ws.Cells.Item(6, 8).Value2 = 143.587201
$ws.Cells.Item(6, 9).Value2 = 0.1228118231805696
$ws.Cells.Item(6, 10).Value2 = 0.1228118231805696
$ws.Cells.Item(6, 13).Value2 = 28.60040066666667
$ws.Cells.Item(6, 14).Value2 = 85.801202
$ws.Cells.Item(6, 15).Value2 = 0.1669535350673389
$ws.Cells.Item(6, 16).Value2 = 0.1669535350673389
$ws.Cells.Item(6, 17).Value2 = 1368.883826401734
$ws.Cells.Item(6, 18).Value2 = 12319.9544376156
$ws.Cells.Item(6, 19).Value2 = 0.02050386802806105
$ws.Cells.Item(6, 20).Value2 = 0.02050386802806106
$ws.Cells.Item(7, 9).Value2 = 0.04786922362394307
$ws.Cells.Item(7, 10).Value2 = 0.04786922362394307
$ws.Cells.Item(7, 13).Value2 = 42.99144133333333
$ws.Cells.Item(7, 14).Value2 = 128.974324
$ws.Cells.Item(7, 15).Value2 = 0.2509605789056467
$ws.Cells.Item(7, 16).Value2 = 0.2509605789056467
$ws.Cells.Item(7, 17).Value2 = 802.033827068457
$ws.Cells.Item(7, 18).Value2 = 7218.304443616112
$ws.Cells.Item(7, 19).Value2 = 0.01201328807242861
$ws.Cells.Item(7, 20).Value2 = 0.01201328807242861
$ws.Cells.Item(8, 9).Value2 = 0.04786922362394307
$ws.Cells.Item(8, 10).Value2 = 0.04786922362394307
$ws.Cells.Item(8, 15).Value2 = 0.3246035785195009
$ws.Cells.Item(8, 16).Value2 = 0.324603578519501
$ws.Cells.Item(8, 19).Value2 = 0.01553852128928216
$ws.Cells.Item(8, 20).Value2 = 0.01553852128928216
$ws.Cells.Item(9, 9).Value2 = 0.04786922362394307
$ws.Cells.Item(9, 10).Value2 = 0.04786922362394307
$ws.Cells.Item(9, 13).Value2 = 36.72715
$ws.Cells.Item(9, 14).Value2 = 110.18145
$ws.Cells.Item(9, 15).Value2 = 0.2143930638214748
$ws.Cells.Item(9, 16).Value2 = 0.2143930638214748
$ws.Cells.Item(9, 17).Value2 = 685.1693211080667
$ws.Cells.Item(9, 18).Value2 = 6166.523889972601
$ws.Cells.Item(9, 19).Value2 = 0.01026282951549247
$ws.Cells.Item(9, 20).Value2 = 0.01026282951549248
$ws.Cells.Item(10, 9).Value2 = 0.04786922362394307
$ws.Cells.Item(10, 10).Value2 = 0.04786922362394307
$ws.Cells.Item(10, 13).Value2 = 7.381512666666667
$ws.Cells.Item(10, 14).Value2 = 22.144538
$ws.Cells.Item(10, 15).Value2 = 0.04308924368603855
$ws.Cells.Item(10, 16).Value2 = 0.04308924368603857
$ws.Cells.Item(10, 17).Value2 = 137.7070102790605
$ws.Cells.Item(10, 18).Value2 = 1239.363092511544
$ws.Cells.Item(10, 19).Value2 = 0.002062648641793556
$ws.Cells.Item(10, 20).Value2 = 0.002062648641793557
$ws.Cells.Item(11, 9).Value2 = 0.04786922362394307
$ws.Cells.Item(11, 10).Value2 = 0.04786922362394307
$ws.Cells.Item(11, 13).Value2 = 28.60040066666667
$ws.Cells.Item(11, 14).Value2 = 85.801202
$ws.Cells.Item(11, 15).Value2 = 0.1669535350673389
$ws.Cells.Item(11, 16).Value2 = 0.1669535350673389
$ws.Cells.Item(11, 17).Value2 = 533.5594269688419
$ws.Cells.Item(11, 18).Value2 = 4802.034842719576
$ws.Cells.Item(11, 19).Value2 = 0.007991936104946267
$ws.Cells.Item(11, 20).Value2 = 0.007991936104946269
$ws.Cells.Item(12, 7).Value2 = 171.0598806666667
$ws.Cells.Item(12, 8).Value2 = 513.1796420000001
$ws.Cells.Item(12, 9).Value2 = 0.4389285884413335
$ws.Cells.Item(12, 10).Value2 = 0.4389285884413335
$ws.Cells.Item(12, 13).Value2 = 42.99144133333333
$ws.Cells.Item(12, 14).Value2 = 128.974324
$ws.Cells.Item(12, 15).Value2 = 0.2509605789056467
$ws.Cells.Item(12, 16).Value2 = 0.2509605789056467
$ws.Cells.Item(12, 17).Value2 = 7354.110824168002
$ws.Cells.Item(12, 18).Value2 = 66186.99741751202
$ws.Cells.Item(12, 19).Value2 = 0.1101537726534754
$ws.Cells.Item(12, 20).Value2 = 0.1101537726534754
$ws.Cells.Item(13, 7).Value2 = 171.0598806666667
$ws.Cells.Item(13, 8).Value2 = 513.1796420000001
$ws.Cells.Item(13, 9).Value2 = 0.4389285884413335
$ws.Cells.Item(13, 10).Value2 = 0.4389285884413335
$ws.Cells.Item(13, 15).Value2 = 0.3246035785195009
$ws.Cells.Item(13, 16).Value2 = 0.324603578519501
$ws.Cells.Item(13, 17).Value2 = 9512.134139806203
$ws.Cells.Item(13, 18).Value2 = 85609.20725825582
$ws.Cells.Item(13, 19).Value2 = 0.1424777905225701
$ws.Cells.Item(13, 20).Value2 = 0.1424777905225701
$ws.Cells.Item(14, 7).Value2 = 171.0598806666667
$ws.Cells.Item(14, 8).Value2 = 513.1796420000001
$ws.Cells.Item(14, 9).Value2 = 0.4389285884413335
$ws.Cells.Item(14, 10).Value2 = 0.4389285884413335
$ws.Cells.Item(14, 13).Value2 = 36.72715
$ws.Cells.Item(14, 14).Value2 = 110.18145
$ws.Cells.Item(14, 15).Value2 = 0.2143930638214748
$ws.Cells.Item(14, 16).Value2 = 0.2143930638214748
$ws.Cells.Item(14, 17).Value2 = 6282.541896226768
$ws.Cells.Item(14, 18).Value2 = 56542.87706604091
$ws.Cells.Item(14, 19).Value2 = 0.09410324487477265
$ws.Cells.Item(14, 20).Value2 = 0.09410324487477266
$ws.Cells.Item(15, 7).Value2 = 171.0598806666667
$ws.Cells.Item(15, 8).Value2 = 513.1796420000001
$ws.Cells.Item(15, 9).Value2 = 0.4389285884413335
$ws.Cells.Item(15, 10).Value2 = 0.4389285884413335
$ws.Cells.Item(15, 13).Value2 = 7.381512666666667
$ws.Cells.Item(15, 14).Value2 = 22.144538
$ws.Cells.Item(15, 15).Value2 = 0.04308924368603855
$ws.Cells.Item(15, 16).Value2 = 0.04308924368603857
$ws.Cells.Item(15, 17).Value2 = 1262.680675899489
$ws.Cells.Item(15, 18).Value2 = 11364.1260830954
$ws.Cells.Item(15, 19).Value2 = 0.01891310090811754
$ws.Cells.Item(15, 20).Value2 = 0.01891310090811755
$ws.Cells.Item(16, 7).Value2 = 171.0598806666667
$ws.Cells.Item(16, 8).Value2 = 513.1796420000001
$ws.Cells.Item(16, 9).Value2 = 0.4389285884413335
$ws.Cells.Item(16, 10).Value2 = 0.4389285884413335
$ws.Cells.Item(16, 13).Value2 = 28.60040066666667
$ws.Cells.Item(16, 14).Value2 = 85.801202
$ws.Cells.Item(16, 15).Value2 = 0.1669535350673389
$ws.Cells.Item(16, 16).Value2 = 0.1669535350673389
$ws.Cells.Item(16, 17).Value2 = 4892.381125058855
$ws.Cells.Item(16, 18).Value2 = 44031.43012552969
$ws.Cells.Item(16, 19).Value2 = 0.07328067948239773
$ws.Cells.Item(16, 20).Value2 = 0.07328067948239775
$ws.Cells.Item(17, 7).Value2 = 12.628047
$ws.Cells.Item(17, 8).Value2 = 37.884141
$ws.Cells.Item(17, 9).Value2 = 0.0324027517316099
$ws.Cells.Item(17, 10).Value2 = 0.0324027517316099
$ws.Cells.Item(17, 13).Value2 = 42.99144133333333
$ws.Cells.Item(17, 14).Value2 = 128.974324
$ws.Cells.Item(17, 15).Value2 = 0.2509605789056467
$ws.Cells.Item(17, 16).Value2 = 0.2509605789056467
$ws.Cells.Item(17, 17).Value2 = 542.897941755076
$ws.Cells.Item(17, 18).Value2 = 4886.081475795684
$ws.Cells.Item(17, 19).Value2 = 0.008131813332700764
$ws.Cells.Item(17, 20).Value2 = 0.008131813332700766
$ws.Cells.Item(18, 7).Value2 = 12.628047
$ws.Cells.Item(18, 8).Value2 = 37.884141
$ws.Cells.Item(18, 9).Value2 = 0.0324027517316099
$ws.Cells.Item(18, 10).Value2 = 0.0324027517316099
$ws.Cells.Item(18, 15).Value2 = 0.3246035785195009
$ws.Cells.Item(18, 16).Value2 = 0.324603578519501
$ws.Cells.Item(18, 17).Value2 = 702.2083525350209
$ws.Cells.Item(18, 18).Value2 = 6319.875172815188
$ws.Cells.Item(18, 19).Value2 = 0.01051804916595953
$ws.Cells.Item(18, 20).Value2 = 0.01051804916595953
$ws.Cells.Item(19, 7).Value2 = 12.628047
$ws.Cells.Item(19, 8).Value2 = 37.884141
$ws.Cells.Item(19, 9).Value2 = 0.0324027517316099
$ws.Cells.Item(19, 10).Value2 = 0.0324027517316099
$ws.Cells.Item(19, 13).Value2 = 36.72715
$ws.Cells.Item(19, 14).Value2 = 110.18145
$ws.Cells.Item(19, 15).Value2 = 0.2143930638214748
$ws.Cells.Item(19, 16).Value2 = 0.2143930638214748
$ws.Cells.Item(19, 17).Value2 = 463.7921763760501
$ws.Cells.Item(19, 18).Value2 = 4174.12958738445
$ws.Cells.Item(19, 19).Value2 = 0.006946925219986443
$ws.Cells.Item(19, 20).Value2 = 0.006946925219986444
$ws.Cells.Item(20, 7).Value2 = 12.628047
$ws.Cells.Item(20, 8).Value2 = 37.884141
$ws.Cells.Item(20, 9).Value2 = 0.0324027517316099
$ws.Cells.Item(20, 10).Value2 = 0.0324027517316099
$ws.Cells.Item(20, 13).Value2 = 7.381512666666667
$ws.Cells.Item(20, 14).Value2 = 22.144538
$ws.Cells.Item(20, 15).Value2 = 0.04308924368603855
$ws.Cells.Item(20, 16).Value2 = 0.04308924368603857
$ws.Cells.Item(20, 17).Value2 = 93.21408888576201
$ws.Cells.Item(20, 18).Value2 = 838.926799971858
$ws.Cells.Item(20, 19).Value2 = 0.001396210065461546
$ws.Cells.Item(20, 20).Value2 = 0.001396210065461547
$ws.Cells.Item(21, 7).Value2 = 12.628047
$ws.Cells.Item(21, 8).Value2 = 37.884141
$ws.Cells.Item(21, 9).Value2 = 0.0324027517316099
$ws.Cells.Item(21, 10).Value2 = 0.0324027517316099
$ws.Cells.Item(21, 13).Value2 = 28.60040066666667
$ws.Cells.Item(21, 14).Value2 = 85.801202
$ws.Cells.Item(21, 15).Value2 = 0.1669535350673389
$ws.Cells.Item(21, 16).Value2 = 0.1669535350673389
$ws.Cells.Item(21, 17).Value2 = 361.1672038374981
$ws.Cells.Item(21, 18).Value2 = 3250.504834537482
$ws.Cells.Item(21, 19).Value2 = 0.005409753947501609
$ws.Cells.Item(21, 20).Value2 = 0.00540975394750161
$ws.Cells.Item(22, 7).Value2 = 139.5154473333333
$ws.Cells.Item(22, 8).Value2 = 418.546342
$ws.Cells.Item(22, 9).Value2 = 0.3579876130225438
$ws.Cells.Item(22, 10).Value2 = 0.3579876130225438
$ws.Cells.Item(22, 13).Value2 = 42.99144133333333
$ws.Cells.Item(22, 14).Value2 = 128.974324
$ws.Cells.Item(22, 15).Value2 = 0.2509605789056467
$ws.Cells.Item(22, 16).Value2 = 0.2509605789056467
$ws.Cells.Item(22, 17).Value2 = 5997.970169124756
$ws.Cells.Item(22, 18).Value2 = 53981.73152212281
$ws.Cells.Item(22, 19).Value2 = 0.08984077860518822
$ws.Cells.Item(22, 20).Value2 = 0.08984077860518823
$ws.Cells.Item(23, 7).Value2 = 139.5154473333333
$ws.Cells.Item(23, 8).Value2 = 418.546342
$ws.Cells.Item(23, 9).Value2 = 0.3579876130225438
$ws.Cells.Item(23, 10).Value2 = 0.3579876130225438
$ws.Cells.Item(23, 15).Value2 = 0.3246035785195009
$ws.Cells.Item(23, 16).Value2 = 0.324603578519501
$ws.Cells.Item(23, 17).Value2 = 7758.041479028901
$ws.Cells.Item(23, 18).Value2 = 69822.3733112601
$ws.Cells.Item(23, 19).Value2 = 0.116204060252772
$ws.Cells.Item(23, 20).Value2 = 0.1162040602527721
$ws.Cells.Item(24, 7).Value2 = 139.5154473333333
$ws.Cells.Item(24, 8).Value2 = 418.546342
$ws.Cells.Item(24, 9).Value2 = 0.3579876130225438
$ws.Cells.Item(24, 10).Value2 = 0.3579876130225438
$ws.Cells.Item(24, 13).Value2 = 36.72715
$ws.Cells.Item(24, 14).Value2 = 110.18145
$ws.Cells.Item(24, 15).Value2 = 0.2143930638214748
$ws.Cells.Item(24, 16).Value2 = 0.2143930638214748
$ws.Cells.Item(24, 17).Value2 = 5124.004761528433
$ws.Cells.Item(24, 18).Value2 = 46116.04285375591
$ws.Cells.Item(24, 19).Value2 = 0.07675006116603965
$ws.Cells.Item(24, 20).Value2 = 0.07675006116603966
$ws.Cells.Item(25, 7).Value2 = 139.5154473333333
$ws.Cells.Item(25, 8).Value2 = 418.546342
$ws.Cells.Item(25, 9).Value2 = 0.3579876130225438
$ws.Cells.Item(25, 10).Value2 = 0.3579876130225438
$ws.Cells.Item(25, 13).Value2 = 7.381512666666667
$ws.Cells.Item(25, 14).Value2 = 22.144538
$ws.Cells.Item(25, 15).Value2 = 0.04308924368603855
$ws.Cells.Item(25, 16).Value2 = 0.04308924368603857
$ws.Cells.Item(25, 17).Value2 = 1029.835041686666
$ws.Cells.Item(25, 18).Value2 = 9268.515375179995
$ws.Cells.Item(25, 19).Value2 = 0.01542541549411166
$ws.Cells.Item(25, 20).Value2 = 0.01542541549411166
$ws.Cells.Item(26, 7).Value2 = 139.5154473333333
$ws.Cells.Item(26, 8).Value2 = 418.546342
$ws.Cells.Item(26, 9).Value2 = 0.3579876130225438
$ws.Cells.Item(26, 10).Value2 = 0.3579876130225438
$ws.Cells.Item(26, 13).Value2 = 28.60040066666667
$ws.Cells.Item(26, 14).Value2 = 85.801202
$ws.Cells.Item(26, 15).Value2 = 0.1669535350673389
$ws.Cells.Item(26, 16).Value2 = 0.1669535350673389
$ws.Cells.Item(26, 17).Value2 = 3990.197692922565
$ws.Cells.Item(26, 18).Value2 = 35911.77923630308
$ws.Cells.Item(26, 19).Value2 = 0.05976729750443222
$ws.Cells.Item(26, 20).Value2 = 0.05976729750443223
